# Apply "Added data for 2 new cities" edit to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows 32-34 (corrected bird_richness / mammal_richness values) ---
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = 4

$ws.Range("E33").Value = 12

$ws.Range("E34").Value = 9

# --- Append 2 new city rows ---
# Row 35: Salvador
$ws.Range("A35").Value = "Salvador"
$ws.Range("B35").Value = 10
$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 15
$ws.Range("F35").Value = 5

# Row 36: Fortaleza
$ws.Range("A36").Value = "Fortaleza"
$ws.Range("B36").Value = 21
$ws.Range("C36").Value = 6
$ws.Range("D36").Value = 10
$ws.Range("E36").Value = 11
$ws.Range("F36").Value = 2

# --- Recolor rows 32:36 (A:F) with a dark-red font (new style, Calibri Light / FFC00000) ---
$range = $ws.Range("A32:F36")
$range.Font.Color = 192

# --- Update selection to match the saved view state ---
$ws.Range("H28").Select()
